# Insert a new "First Paragraph"-styled paragraph describing the optional
# "useBom" / "pushAfterConfigChange" configuration keys. It is placed right
# after the "База данни (изисква IT специалист)" bullet and right before the
# "Използване на Електронна таблица (Microsoft Excel документ)" Heading3.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "База данни (изисква IT специалист)") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'База данни (изисква IT специалист)' paragraph"
}

# Create a new, empty paragraph right after $target.
$tr = $target.Range
$tr.Collapse(0)
$tr.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Style = "First Paragraph"
$nr = $newPara.Range

# Appends $text to the end of $range's paragraph content (i.e. right before
# the trailing paragraph mark), optionally making just that chunk bold.
function Add-Run($range, [string]$text, [bool]$bold) {
    $startPos = $range.End - 1
    $range.InsertAfter($text)
    if ($bold) {
        $endPos = $range.End - 1
        $chunk = $d.Range($startPos, $endPos)
        $chunk.Font.Bold = 1
    }
}

Add-Run $nr "Освен конфигурационните опции посочени по-горе, има и още две незадължителни:" $false
Add-Run $nr " " $false
Add-Run $nr "useBom" $true
Add-Run $nr " " $false
Add-Run $nr "и" $false
Add-Run $nr " " $false
Add-Run $nr "pushAfterConfigChange" $true
Add-Run $nr ". Те" $false
Add-Run $nr " " $false
Add-Run $nr "приемат стойности" $false
Add-Run $nr " " $false
Add-Run $nr "true" $true
Add-Run $nr " " $false
Add-Run $nr "или" $false
Add-Run $nr " " $false
Add-Run $nr "false" $true
Add-Run $nr " " $false
Add-Run $nr "и се отнасят съответно за добавянето на UTF-8 BOM във файловете, и дали даден" $false
Add-Run $nr " " $false
Add-Run $nr "ресурс да бъде публикуван при промяна в конфигурацията. По подразбиране стойността е" $false
Add-Run $nr " " $false
Add-Run $nr "false" $true
Add-Run $nr "." $false
